# "Quick message testcases module"
# Flip the enabled/disabled flags in the "Suite" sheet:
#   - Web_QUICK_MESSAGES (row 4, col C) : N -> Y  (now included in the suite)
#   - Web_PATIENT        (row 6, col C) : Y -> N  (no longer included)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Y"
$ws.Range("C6").Value = "N"

# Matches the saved cursor/selection position recorded in the edited file.
$ws.Range("D23").Select()
